$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (Beta) - new unfolding results with 100 keV threshold
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 7.070407707669703
$ws.Range("G2").Value = 3.175127501275629
$ws.Range("H2").Value = 10.58953095301975
$ws.Range("I2").Value = 0.3025104915278049
$ws.Range("J2").Value = 0.2492545487775952
$ws.Range("K2").Value = 0.3438127199915785
$ws.Range("L2").Value = 0.04320990133013451
$ws.Range("M2").Value = 0.0320850825649982
$ws.Range("N2").Value = 0.05349341376574289

# Update existing row 3 values (Gamma)
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.2455320766091476
$ws.Range("G3").Value = 0.001619681515902877
$ws.Range("H3").Value = 0.575183570330392
$ws.Range("I3").Value = 0.2275789218708711
$ws.Range("J3").Value = 0.001491703297717326
$ws.Range("K3").Value = 0.5337275622022035
$ws.Range("L3").Value = 0.2562564349606977
$ws.Range("M3").Value = 0.001715461819326484
$ws.Range("N3").Value = 0.5986381997956087

# Add new row 4 (Beta + Gamma) with 100 keV threshold results
$ws.Range("A4").Value = 2
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 7.31593978427885
$ws.Range("G4").Value = 3.176747182791532
$ws.Range("H4").Value = 11.16471452335015
$ws.Range("I4").Value = 0.530089413398676
$ws.Range("J4").Value = 0.2507462520753125
$ws.Range("K4").Value = 0.877540282193782
$ws.Range("L4").Value = 0.2994663362908322
$ws.Range("M4").Value = 0.03380054438432468
$ws.Range("N4").Value = 0.6521316135613515
